# edit.ps1
# Applies the betexplorer re-scrape update for
# saudi-arabia_saudi-professional-league_2023-2024:
#   - re-sync reordered a few already-scraped fixtures (same matches,
#     rows 21-24 / 83-86 / 92-93 now carry a different fixtures odds
#     in columns F:V; columns A:E -- the row index/league/date key --
#     are untouched)
#   - 3 new fixtures were appended as rows 101-103

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-synced rows: only F:V change, A:E (index/country/league/season/date) stay put ---
# Row 21: was row 24's fixture (Al Ahli SC vs Al Akhdoud)
$ws.Range("F21").Value = "Al Ahli SC"
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = "Al Akhdoud"
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1.16
$ws.Range("K21").Value = "22/08/2023 07:46"
$ws.Range("L21").Value = 1.11
$ws.Range("M21").Value = "24/08/2023 19:40"
$ws.Range("N21").Value = 7.91
$ws.Range("O21").Value = "22/08/2023 07:46"
$ws.Range("P21").Value = 9.279999999999999
$ws.Range("Q21").Value = "24/08/2023 19:52"
$ws.Range("R21").Value = 16.84
$ws.Range("S21").Value = "22/08/2023 07:46"
$ws.Range("T21").Value = 29.41
$ws.Range("U21").Value = "24/08/2023 19:52"
$ws.Range("V21").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ahli-sc-al-akhdoud/IgjeCGQr/"

# Row 22: was row 21's fixture (Al Riyadh vs Al Ittihad)
$ws.Range("F22").Value = "Al Riyadh"
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = "Al Ittihad"
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 13.51
$ws.Range("K22").Value = "22/08/2023 07:46"
$ws.Range("L22").Value = 12.1
$ws.Range("M22").Value = "24/08/2023 19:54"
$ws.Range("N22").Value = 6.52
$ws.Range("O22").Value = "22/08/2023 07:46"
$ws.Range("P22").Value = 6.94
$ws.Range("Q22").Value = "24/08/2023 19:54"
$ws.Range("R22").Value = 1.21
$ws.Range("S22").Value = "22/08/2023 07:46"
$ws.Range("T22").Value = 1.22
$ws.Range("U22").Value = "24/08/2023 19:07"
$ws.Range("V22").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-riyadh-al-ittihad/zVp0Bztk/"

# Row 23: was row 22's fixture (Al Raed vs Al Hilal)
$ws.Range("F23").Value = "Al Raed"
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = "Al Hilal"
$ws.Range("I23").Value = 4
$ws.Range("J23").Value = 12.64
$ws.Range("K23").Value = "22/08/2023 07:46"
$ws.Range("L23").Value = 18.25
$ws.Range("M23").Value = "24/08/2023 19:59"
$ws.Range("N23").Value = 6.75
$ws.Range("O23").Value = "22/08/2023 07:46"
$ws.Range("P23").Value = 9.07
$ws.Range("Q23").Value = "24/08/2023 19:59"
$ws.Range("R23").Value = 1.21
$ws.Range("S23").Value = "22/08/2023 07:46"
$ws.Range("T23").Value = 1.14
$ws.Range("U23").Value = "24/08/2023 19:51"
$ws.Range("V23").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-raed-al-hilal/MN4PHx3L/"

# Row 24: was row 23's fixture (Al Ettifaq vs Al Khaleej)
$ws.Range("F24").Value = "Al Ettifaq"
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = "Al Khaleej"
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 1.56
$ws.Range("K24").Value = "22/08/2023 07:46"
$ws.Range("L24").Value = 1.85
$ws.Range("M24").Value = "24/08/2023 19:54"
$ws.Range("N24").Value = 4.25
$ws.Range("O24").Value = "22/08/2023 07:46"
$ws.Range("P24").Value = 3.7
$ws.Range("Q24").Value = "24/08/2023 19:54"
$ws.Range("R24").Value = 4.92
$ws.Range("S24").Value = "22/08/2023 07:46"
$ws.Range("T24").Value = 4.2
$ws.Range("U24").Value = "24/08/2023 19:54"
$ws.Range("V24").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-ettifaq-fc-al-khaleej/Gp4TGdIR/"

# Row 83: was row 84's fixture (Al Taawon vs Al Ittihad)
$ws.Range("F83").Value = "Al Taawon"
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = "Al Ittihad"
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 4.56
$ws.Range("K83").Value = "15/10/2023 12:47"
$ws.Range("L83").Value = 4.61
$ws.Range("M83").Value = "20/10/2023 16:58"
$ws.Range("N83").Value = 4.09
$ws.Range("O83").Value = "15/10/2023 12:47"
$ws.Range("P83").Value = 4.24
$ws.Range("Q83").Value = "20/10/2023 16:58"
$ws.Range("R83").Value = 1.63
$ws.Range("S83").Value = "15/10/2023 12:47"
$ws.Range("T83").Value = 1.68
$ws.Range("U83").Value = "20/10/2023 16:51"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taawon-al-ittihad/44sRvIN0/"

# Row 84: was row 83's fixture (Al Akhdoud vs Al Feiha)
$ws.Range("F84").Value = "Al Akhdoud"
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = "Al Feiha"
$ws.Range("I84").Value = 2
$ws.Range("J84").Value = 3.22
$ws.Range("K84").Value = "15/10/2023 17:20"
$ws.Range("L84").Value = 2.86
$ws.Range("M84").Value = "20/10/2023 16:58"
$ws.Range("N84").Value = 3.32
$ws.Range("O84").Value = "15/10/2023 17:20"
$ws.Range("P84").Value = 3.49
$ws.Range("Q84").Value = "20/10/2023 16:59"
$ws.Range("R84").Value = 2.25
$ws.Range("S84").Value = "15/10/2023 17:20"
$ws.Range("T84").Value = 2.45
$ws.Range("U84").Value = "20/10/2023 16:59"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-akhdoud-al-feiha/EZOH1uVD/"

# Row 85: was row 86's fixture (Al Fateh vs Abha)
$ws.Range("F85").Value = "Al Fateh"
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = "Abha"
$ws.Range("I85").Value = 1
$ws.Range("J85").Value = 1.54
$ws.Range("K85").Value = "15/10/2023 12:47"
$ws.Range("L85").Value = 1.4
$ws.Range("M85").Value = "20/10/2023 19:50"
$ws.Range("N85").Value = 4.56
$ws.Range("O85").Value = "15/10/2023 12:47"
$ws.Range("P85").Value = 5.34
$ws.Range("Q85").Value = "20/10/2023 19:50"
$ws.Range("R85").Value = 4.75
$ws.Range("S85").Value = "15/10/2023 12:47"
$ws.Range("T85").Value = 6.72
$ws.Range("U85").Value = "20/10/2023 19:50"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-fateh-abha/bRtVwxw7/"

# Row 86: was row 85's fixture (Al Hilal vs Al Khaleej)
$ws.Range("F86").Value = "Al Hilal"
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = "Al Khaleej"
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1.1
$ws.Range("K86").Value = "15/10/2023 12:47"
$ws.Range("L86").Value = 1.1
$ws.Range("M86").Value = "20/10/2023 19:38"
$ws.Range("N86").Value = 11.38
$ws.Range("O86").Value = "15/10/2023 12:47"
$ws.Range("P86").Value = 11.29
$ws.Range("Q86").Value = "20/10/2023 19:53"
$ws.Range("R86").Value = 21.56
$ws.Range("S86").Value = "15/10/2023 12:47"
$ws.Range("T86").Value = 19.58
$ws.Range("U86").Value = "20/10/2023 19:53"
$ws.Range("V86").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-hilal-al-khaleej/0jNPaJaQ/"

# Row 92: was row 93's fixture (Damac vs Al Akhdoud)
$ws.Range("F92").Value = "Damac"
$ws.Range("G92").Value = 2
$ws.Range("H92").Value = "Al Akhdoud"
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1.88
$ws.Range("K92").Value = "24/10/2023 21:01"
$ws.Range("L92").Value = 1.84
$ws.Range("M92").Value = "26/10/2023 16:50"
$ws.Range("N92").Value = 3.79
$ws.Range("O92").Value = "24/10/2023 21:01"
$ws.Range("P92").Value = 3.95
$ws.Range("Q92").Value = "26/10/2023 17:00"
$ws.Range("R92").Value = 3.92
$ws.Range("S92").Value = "24/10/2023 21:01"
$ws.Range("T92").Value = 3.99
$ws.Range("U92").Value = "26/10/2023 16:50"
$ws.Range("V92").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/damac-al-akhdoud/S6fHkZyJ/"

# Row 93: was row 92's fixture (Al Khaleej vs Al Taawon)
$ws.Range("F93").Value = "Al Khaleej"
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = "Al Taawon"
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = 3.53
$ws.Range("K93").Value = "24/10/2023 21:01"
$ws.Range("L93").Value = 3.5
$ws.Range("M93").Value = "26/10/2023 16:55"
$ws.Range("N93").Value = 3.77
$ws.Range("O93").Value = "24/10/2023 21:01"
$ws.Range("P93").Value = 3.7
$ws.Range("Q93").Value = "26/10/2023 16:56"
$ws.Range("R93").Value = 2
$ws.Range("S93").Value = "24/10/2023 21:01"
$ws.Range("T93").Value = 2.04
$ws.Range("U93").Value = "26/10/2023 16:55"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-khaleej-al-taawon/Wv18ie76/"

# --- New fixtures appended at the bottom of the table ---
# Row 101: Al Fateh vs Al Hilal
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = "saudi-arabia"
$ws.Range("C101").Value = "saudi-professional-league"
$ws.Range("D101").Value = "2023-2024"
$ws.Range("E101").Value = 45233.66666666666
$ws.Range("F101").Value = "Al Fateh"
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = "Al Hilal"
$ws.Range("I101").Value = 2
$ws.Range("J101").Value = 5.32
$ws.Range("K101").Value = "29/10/2023 19:43"
$ws.Range("L101").Value = 10.37
$ws.Range("M101").Value = "03/11/2023 15:57"
$ws.Range("N101").Value = 5.18
$ws.Range("O101").Value = "29/10/2023 19:43"
$ws.Range("P101").Value = 6.9
$ws.Range("Q101").Value = "03/11/2023 15:57"
$ws.Range("R101").Value = 1.43
$ws.Range("S101").Value = "29/10/2023 19:43"
$ws.Range("T101").Value = 1.23
$ws.Range("U101").Value = "03/11/2023 15:52"
$ws.Range("V101").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-fateh-al-hilal/0YXc2CLO/"

# Row 102: Al Taee vs Al Feiha
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = "saudi-arabia"
$ws.Range("C102").Value = "saudi-professional-league"
$ws.Range("D102").Value = "2023-2024"
$ws.Range("E102").Value = 45233.66666666666
$ws.Range("F102").Value = "Al Taee"
$ws.Range("G102").Value = 3
$ws.Range("H102").Value = "Al Feiha"
$ws.Range("I102").Value = 3
$ws.Range("J102").Value = 2.73
$ws.Range("K102").Value = "29/10/2023 19:43"
$ws.Range("L102").Value = 3.06
$ws.Range("M102").Value = "03/11/2023 15:55"
$ws.Range("N102").Value = 3.38
$ws.Range("O102").Value = "29/10/2023 19:43"
$ws.Range("P102").Value = 3.43
$ws.Range("Q102").Value = "03/11/2023 15:55"
$ws.Range("R102").Value = 2.45
$ws.Range("S102").Value = "29/10/2023 19:43"
$ws.Range("T102").Value = 2.35
$ws.Range("U102").Value = "03/11/2023 15:55"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-taee-al-feiha/SKZk4YjC/"

# Row 103: Al Shabab vs Al Ittihad
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = "saudi-arabia"
$ws.Range("C103").Value = "saudi-professional-league"
$ws.Range("D103").Value = "2023-2024"
$ws.Range("E103").Value = 45233.79166666666
$ws.Range("F103").Value = "Al Shabab"
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = "Al Ittihad"
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 3.74
$ws.Range("K103").Value = "29/10/2023 19:43"
$ws.Range("L103").Value = 3.48
$ws.Range("M103").Value = "03/11/2023 18:50"
$ws.Range("N103").Value = 3.87
$ws.Range("O103").Value = "29/10/2023 19:43"
$ws.Range("P103").Value = 3.52
$ws.Range("Q103").Value = "03/11/2023 18:59"
$ws.Range("R103").Value = 1.91
$ws.Range("S103").Value = "29/10/2023 19:43"
$ws.Range("T103").Value = 1.95
$ws.Range("U103").Value = "03/11/2023 18:50"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/saudi-arabia/saudi-professional-league/al-shabab-al-ittihad/2ZP0MF6t/"

# Formatting for the new rows, matching the existing table conventions:
$colA = $ws.Range("A101:A103")
$colA.Font.Bold = $true
$colA.HorizontalAlignment = -4108
$colA.VerticalAlignment = -4160
$colA.Borders.LineStyle = 1
$ws.Range("E101:E103").NumberFormat = "YYYY-MM-DD HH:MM:SS"

